# Update the two-digit / one-digit division worksheet values in the
# single table of the document. Each data row of the table sits at
# table-row indices 1, 5, 9, 13, 17 (the intervening rows are blank
# spacer rows). We assign each target cells Range.Text directly
# (rather than using Find/Replace) so that a new value which happens
# to equal another cells old value cannot be matched twice.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "94÷6=") {
    $cell.Range.Text = "75÷6="
}
$cell = $t.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "57÷2=") {
    $cell.Range.Text = "42÷4="
}
$cell = $t.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "32÷7=") {
    $cell.Range.Text = "58÷3="
}
$cell = $t.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "33÷8=") {
    $cell.Range.Text = "50÷2="
}
$cell = $t.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "46÷7=") {
    $cell.Range.Text = "44÷9="
}

$cell = $t.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "91÷5=") {
    $cell.Range.Text = "26÷2="
}
$cell = $t.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "88÷3=") {
    $cell.Range.Text = "47÷8="
}
$cell = $t.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "63÷5=") {
    $cell.Range.Text = "88÷6="
}
$cell = $t.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "50÷4=") {
    $cell.Range.Text = "42÷5="
}
$cell = $t.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "51÷2=") {
    $cell.Range.Text = "53÷6="
}

$cell = $t.Cell(9, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "61÷2=") {
    $cell.Range.Text = "23÷9="
}
$cell = $t.Cell(9, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "21÷2=") {
    $cell.Range.Text = "63÷7="
}
$cell = $t.Cell(9, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "93÷8=") {
    $cell.Range.Text = "19÷3="
}
$cell = $t.Cell(9, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "70÷4=") {
    $cell.Range.Text = "62÷6="
}
$cell = $t.Cell(9, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "86÷5=") {
    $cell.Range.Text = "74÷8="
}

$cell = $t.Cell(13, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "46÷3=") {
    $cell.Range.Text = "61÷3="
}
$cell = $t.Cell(13, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "62÷4=") {
    $cell.Range.Text = "50÷4="
}
$cell = $t.Cell(13, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "54÷8=") {
    $cell.Range.Text = "45÷7="
}
$cell = $t.Cell(13, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "84÷5=") {
    $cell.Range.Text = "24÷7="
}
$cell = $t.Cell(13, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "17÷2=") {
    $cell.Range.Text = "31÷7="
}

$cell = $t.Cell(17, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "60÷6=") {
    $cell.Range.Text = "34÷2="
}
$cell = $t.Cell(17, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "66÷6=") {
    $cell.Range.Text = "99÷9="
}
$cell = $t.Cell(17, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "15÷7=") {
    $cell.Range.Text = "87÷7="
}
$cell = $t.Cell(17, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "19÷3=") {
    $cell.Range.Text = "90÷2="
}
$cell = $t.Cell(17, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "33÷7=") {
    $cell.Range.Text = "34÷7="
}

